$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace comma-decimal text values in the SIZE (cm) column (K) with true
# numeric values using a period as the decimal separator.
$ws.Range("K4").Value = 2.5
$ws.Range("K5").Value = 2.5
$ws.Range("K6").Value = 2.5
$ws.Range("K7").Value = 2.5
$ws.Range("K11").Value = 2.5
$ws.Range("K14").Value = 2.5
$ws.Range("K17").Value = 1.5
$ws.Range("K18").Value = 1.5
$ws.Range("K31").Value = 4.5
$ws.Range("K32").Value = 1.5
$ws.Range("K38").Value = 1.5
$ws.Range("K44").Value = 2.5
$ws.Range("K46").Value = 1.5
$ws.Range("K48").Value = 1.5
$ws.Range("K49").Value = 2.5
$ws.Range("K50").Value = 2.5
$ws.Range("K51").Value = 1.5
$ws.Range("K53").Value = 1.5
$ws.Range("K54").Value = 2.5
$ws.Range("K60").Value = 1.5
$ws.Range("K68").Value = 1.5
$ws.Range("K69").Value = 2.5
$ws.Range("K72").Value = 4.5
$ws.Range("K74").Value = 1.5
$ws.Range("K82").Value = 1.5
$ws.Range("K91").Value = 1.5

# Mirror the end-user's final selection state: column K selected, view
# scrolled back to the top of the sheet.
$ws.Columns("K:K").Select()
